# Add two new paragraphs (plus a blank separator paragraph) right after the
# paragraph ending in "...最长前缀匹配来确定应该匹配哪一个。" and right before the
# existing (empty) "_GoBack" bookmark paragraph that follows it.

$d = $word.ActiveDocument

$anchorText = "最长前缀匹配来确定应该匹配哪一个"

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$anchorText*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the anchor paragraph containing '$anchorText'."
}

# The paragraph right after the anchor is where the new content must be
# inserted in front of, so the new paragraphs land between them.
$insertBeforePara = $d.Paragraphs.Item($targetIndex + 1)
$insertRange = $insertBeforePara.Range

$blockText = "IP 地址和 MAC 地址"
$blockText2 = "网络层实现主机之间的通信，而链路层实现具体每段链路之间的通信。因此在通信过程中，IP 数据报的源地址和目的地址始终不变，而 MAC 地址随着链路的改变而改变。"

# One blank paragraph, then the short "IP address / MAC address" heading-like
# paragraph, then the explanatory paragraph -- each terminated with a
# paragraph mark so they land as separate <w:p> elements ahead of the
# existing following paragraph.
$fullInsert = "`r" + $blockText + "`r" + $blockText2 + "`r"

$insertRange.InsertBefore($fullInsert)
